$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F34").Value = -3.161960781933703
$ws.Range("F35").Value = -3.166162632831467
$ws.Range("F38").Value = -3.178121789665747
$ws.Range("F39").Value = -3.184380113057521
$ws.Range("F42").Value = -3.196985665750812
$ws.Range("F43").Value = -3.201187516648575
$ws.Range("F44").Value = -3.205389367546339
$ws.Range("F48").Value = -3.216556464597398
$ws.Range("F49").Value = -3.190575958804059
$ws.Range("F50").Value = -3.187103221535859
$ws.Range("F53").Value = -3.180699468474087
$ws.Range("F57").Value = -3.155419697303079
$ws.Range("F81").Value = -3.320804007335938
$ws.Range("F82").Value = -3.361833447644869
$ws.Range("F83").Value = -3.376474829882223
$ws.Range("F85").Value = -3.3370655349794
$ws.Range("F87").Value = -3.241838972347475
$ws.Range("F88").Value = -3.166608013075078
$ws.Range("F93").Value = -3.43911252864328
$ws.Range("F99").Value = -3.278502416014376
$ws.Range("F107").Value = -3.485295511903248
$ws.Range("F108").Value = -3.44478978158482
$ws.Range("F110").Value = -3.330813995624236
$ws.Range("F111").Value = -3.255608016254619
$ws.Range("F118").Value = -3.518918828605665
$ws.Range("F120").Value = -3.434275912249621
$ws.Range("F129").Value = -3.552251047350676
$ws.Range("F137").Value = -3.604893142865736
$ws.Range("F138").Value = -3.607889149924412
$ws.Range("F139").Value = -3.596806518493822
$ws.Range("F140").Value = -3.569184181595075
$ws.Range("F141").Value = -3.536419451436262
$ws.Range("F142").Value = -3.489150396009545
$ws.Range("F144").Value = -3.364608870239705
$ws.Range("F147").Value = -3.617589867200584
$ws.Range("F148").Value = -3.624655280125956
$ws.Range("F151").Value = -3.590136156912799
$ws.Range("F154").Value = -3.453243361829357
$ws.Range("F158").Value = -3.632600993004166
$ws.Range("F159").Value = -3.642638051721088
$ws.Range("F160").Value = -3.639368617781685
$ws.Range("F161").Value = -3.629048546142379
$ws.Range("F162").Value = -3.603530931925687
$ws.Range("F166").Value = -3.400228952407928
$ws.Range("F169").Value = -3.648616127823166
$ws.Range("F172").Value = -3.643521795797566
$ws.Range("F175").Value = -3.538632523183368
$ws.Range("F180").Value = -3.662109141693795
$ws.Range("F182").Value = -3.664247735708846
$ws.Range("F187").Value = -3.492132144394738
$ws.Range("F191").Value = -3.672940893973359
$ws.Range("F193").Value = -3.675257212303423
$ws.Range("F197").Value = -3.557315142814651
$ws.Range("F199").Value = -3.436488095947585
$ws.Range("F202").Value = -3.682169032005127
$ws.Range("F203").Value = -3.687827829873812
$ws.Range("F204").Value = -3.684077120514287
$ws.Range("F206").Value = -3.643377545314912
$ws.Range("F207").Value = -3.609640251107101
$ws.Range("F208").Value = -3.5645668228519
$ws.Range("F210").Value = -3.44543562399414
